$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values would otherwise be auto-parsed as numbers by Excel;
# force them to keep their original plain-text representation (matches source data which
# stores prices as literal text, e.g. "25.895.75" alongside "1.003").
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D35", "D38", "D39", "D41", "D42", "D44", "D45", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.895.75"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.637.45"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "214.51"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "0.5067"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2551"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").Value = "0.06367"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "19.43"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "0.07745"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.270"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.648.18"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "0.5431"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "0.0₅7807"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "64.10"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "25.923.57"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "196.19"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").Value = "4.449"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "9.905"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "1.884"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "141.03"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "0.1187"
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("D27").Value = "6.851"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "15.68"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "0.04930"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "3.248"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "3.170"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "0.8931"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").Value = "1.132.71"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "0.5413"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "0.01551"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "2.547"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "5.573"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("D44").Value = "0.8129"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "99.39"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "1.774.81"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").Value = "54.74"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +0.04%  "
